# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Palta, Feria Lagunitas de Puerto Montt)
# just before the existing row 200, pushing the old rows 200-228 down to
# 202-230.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 200 (shifts 200:228 -> 202:230)
$ws.Rows.Item(200).Resize(2).Insert()

# --- New row 200: "1a nueva(o)" ---
$ws.Cells.Item(200, 1).Value  = 4
$ws.Cells.Item(200, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(200, 3).Value  = "Los Lagos"
$ws.Cells.Item(200, 4).Value  = 44491
$ws.Cells.Item(200, 5).Value  = 10
$ws.Cells.Item(200, 6).Value  = "Fruta"
$ws.Cells.Item(200, 7).Value  = 100106
$ws.Cells.Item(200, 8).Value  = "Oleaginosos"
$ws.Cells.Item(200, 9).Value  = 100106002
$ws.Cells.Item(200, 10).Value = "Palta"
$ws.Cells.Item(200, 11).Value = "Hass"
$ws.Cells.Item(200, 12).Value = "1a nueva(o)"
$ws.Cells.Item(200, 13).Value = 300
$ws.Cells.Item(200, 14).Value = 3900
$ws.Cells.Item(200, 15).Value = 4000
$ws.Cells.Item(200, 16).Value = 3950
$ws.Cells.Item(200, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(200, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(200, 19).Value = 3950
$ws.Cells.Item(200, 20).Value = 1

# --- New row 201: "2a nueva(o)" ---
$ws.Cells.Item(201, 1).Value  = 4
$ws.Cells.Item(201, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(201, 3).Value  = "Los Lagos"
$ws.Cells.Item(201, 4).Value  = 44491
$ws.Cells.Item(201, 5).Value  = 10
$ws.Cells.Item(201, 6).Value  = "Fruta"
$ws.Cells.Item(201, 7).Value  = 100106
$ws.Cells.Item(201, 8).Value  = "Oleaginosos"
$ws.Cells.Item(201, 9).Value  = 100106002
$ws.Cells.Item(201, 10).Value = "Palta"
$ws.Cells.Item(201, 11).Value = "Hass"
$ws.Cells.Item(201, 12).Value = "2a nueva(o)"
$ws.Cells.Item(201, 13).Value = 150
$ws.Cells.Item(201, 14).Value = 3500
$ws.Cells.Item(201, 15).Value = 3500
$ws.Cells.Item(201, 16).Value = 3500
$ws.Cells.Item(201, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(201, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(201, 19).Value = 3500
$ws.Cells.Item(201, 20).Value = 1
